$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text looks like a plain decimal number (e.g. "210.76")
# are pre-formatted as Text so Excel keeps them as literal strings instead of
# silently converting them to numeric values (the source cells are inline
# strings, not numbers).

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.235.69'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.562.89'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '210.76'
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '0.0594'
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").Value = '0.0872'
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("D12").Value = '1.786.19'
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = '1.564.02'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '3.76'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").Value = '27.244.85'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").Value = '61.78'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '217.83'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").Value = '0.0₃0702'
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").Value = '9.39'
$ws.Range("E23").Value = '  +1.66%  '
$ws.Range("D24").Value = '1.94'
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("D25").Value = '151.40'
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").Value = '6.62'
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").Value = '15.00'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("E30").Value = '  +1.78%  '
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").Value = '1.458.18'
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("D34").Value = '3.17'
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("E35").Value = '  +4.70%  '
$ws.Range("D36").Value = '1.62'
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").Value = '0.539'
$ws.Range("E39").Value = '  +1.15%  '
$ws.Range("D40").Value = '5.85'
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("D41").Value = '0.813'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").Value = '2.35'
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("D44").Value = '0.975'
$ws.Range("E44").Value = '  -2.77%  '
$ws.Range("D45").Value = '64.31'
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("D46").Value = '1.75'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").Value = '1.700.46'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").Value = '85.77'
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("D51").Value = '0.0946'
$ws.Range("E51").Value = '  -1.59%  '
